$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.72
$ws.Range("K2").Value = 4.7
$ws.Range("O2").Value = 1.14
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 1.84
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.76
$ws.Range("Y2").Value = 32
$ws.Range("AB2").Value = 16
$ws.Range("AD2").Value = 19.5
$ws.Range("AE2").Value = 46
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 19.5
$ws.Range("AL2").Value = 22
$ws.Range("AN2").Value = 6.2
# Row 3
$ws.Range("H3").Value = 4.1
$ws.Range("N3").Value = 5.2
$ws.Range("P3").Value = 2.42
$ws.Range("R3").Value = 1.58
$ws.Range("S3").Value = 2.64
$ws.Range("U3").Value = 2.44
$ws.Range("X3").Value = 20
$ws.Range("AN3").Value = 9.199999999999999
$ws.Range("AO3").Value = 34
# Row 4
$ws.Range("F4").Value = 1.78
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 4.8
$ws.Range("K4").Value = 4.5
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.49
$ws.Range("U4").Value = 2.24
$ws.Range("V4").Value = 1.26
# Row 5
$ws.Range("G5").Value = 2.94
$ws.Range("H5").Value = 2.96
$ws.Range("K5").Value = 4.8
$ws.Range("P5").Value = 1.53
$ws.Range("Q5").Value = 2.08
$ws.Range("W5").Value = 1.51
# Row 6
$ws.Range("H6").Value = 5.3
$ws.Range("K6").Value = 7.8
$ws.Range("N6").Value = 2.72
$ws.Range("O6").Value = 1.06
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 1.84
$ws.Range("R6").Value = 1.22
# Row 7
$ws.Range("F7").Value = 13.5
$ws.Range("G7").Value = 17
$ws.Range("H7").Value = 1.22
$ws.Range("I7").Value = 1.28
$ws.Range("J7").Value = 6.8
$ws.Range("K7").Value = 8.4
# Row 9
$ws.Range("F9").Value = 2.24
$ws.Range("G9").Value = 2.26
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 3.75
$ws.Range("P9").Value = 1.83
$ws.Range("U9").Value = 2.06
$ws.Range("V9").Value = 1.36
$ws.Range("W9").Value = 1.79
$ws.Range("AE9").Value = 46
$ws.Range("AH9").Value = 18
$ws.Range("AK9").Value = 24
$ws.Range("AO9").Value = 50
# Row 10
$ws.Range("T10").Value = 1.77
$ws.Range("U10").Value = 2.24
$ws.Range("X10").Value = 24
$ws.Range("Y10").Value = 26
$ws.Range("AD10").Value = 23
$ws.Range("AH10").Value = 19.5
$ws.Range("AL10").Value = 27
$ws.Range("AN10").Value = 6.6
$ws.Range("AO10").Value = 75
# Row 11
$ws.Range("O11").Value = 1.18
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 2.32
# Row 12
$ws.Range("L12").Value = 1.4
$ws.Range("N12").Value = 4
$ws.Range("Q12").Value = 1.99
$ws.Range("T12").Value = 1.79
$ws.Range("X12").Value = 15
$ws.Range("AE12").Value = 38
# Row 13
$ws.Range("H13").Value = 1.42
$ws.Range("I13").Value = 1.43
$ws.Range("J13").Value = 5.3
$ws.Range("K13").Value = 5.4
$ws.Range("P13").Value = 2.38
$ws.Range("Q13").Value = 1.71
$ws.Range("V13").Value = 3.3
$ws.Range("Z13").Value = 8.4
$ws.Range("AB13").Value = 30
$ws.Range("AJ13").Value = 280
$ws.Range("AK13").Value = 130
$ws.Range("AN13").Value = 150
